# Adds a new "id" column (F) to the worksheet, populated with a simple
# 1..180 row counter, and updates the selection to cover the full table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("F1").Value = "id"

# Determine how many data rows exist (rows 2..181 in the original data)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 181 }

# Fill F2:F<lastRow> with a sequential id (1, 2, 3, ...)
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = $r - 1
}

# Widen column E a bit (closest attainable width to the authored value)
$ws.Columns.Item(5).ColumnWidth = 19.5

# Reflect the new table bounds in the selection, like the authored workbook
[void]$ws.Range("A2:F$lastRow").Select()

# Best-effort: scroll the view down towards the bottom of the table
$win = $excel.ActiveWindow
$win.ScrollRow = 160
$win.ScrollColumn = 1
